$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data rows (row number -> A, B(label), C, D, E)
$data = @(
    @(2,  0, "line1", 7,  9,  $false),
    @(3,  1, "line2", 9,  8,  $true),
    @(4,  2, "line3", 8,  10, $true),
    @(5,  3, "line4", 8,  11, $true),
    @(6,  4, "line5", 10, 5,  $true),
    @(7,  5, "line6", 12, 8,  $true),
    @(8,  6, "line7", 14, 11, $true),
    @(9,  7, "line8", 16, 9,  $true),
    @(10, 8, "extr1", 5,  12, $true),
    @(11, 9, "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $true),
    @(13, 11, "extr4", 7,  8,  $true),
    @(14, 12, "extr5", 9,  11, $false),
    @(15, 13, "extr6", 7,  11, $false),
    @(16, 14, "extr7", 5,  7,  $true),
    @(17, 15, "extr8", 8,  5,  $true)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# New rows 16 and 17 need the same style as column A in the other data rows
# (bold font, thin border box, centered horizontal, top vertical alignment)
$ws.Range("A16:A17").Font.Bold = $true
$ws.Range("A16:A17").HorizontalAlignment = -4108
$ws.Range("A16:A17").VerticalAlignment = -4160
$ws.Range("A16:A17").Borders.LineStyle = 1
